$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.522.81"
$ws.Range("E2").Value = "  +4.14%  "
$ws.Range("D3").Value = "2.632.49"
$ws.Range("E3").Value = "  +4.17%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.21"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.07"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "2.630.37"
$ws.Range("E9").Value = "  +4.08%  "
$ws.Range("E10").Value = "  +14.26%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.01"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "3.114.19"
$ws.Range("E14").Value = "  +4.18%  "
$ws.Range("E15").Value = "  +10.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.81"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "71.409.73"
$ws.Range("D18").Value = "2.623.83"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "382.63"
$ws.Range("E19").Value = "  +8.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.93"
$ws.Range("E20").Value = "  +5.67%  "
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.80"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("E24").Value = "  +16.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.49"
$ws.Range("E25").Value = "  +5.92%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.76"
$ws.Range("E27").Value = "  +8.75%  "
$ws.Range("D28").Value = "2.764.81"
$ws.Range("E28").Value = "  +2.75%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "0.0₃0964"
$ws.Range("E30").Value = "  +7.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "549.06"
$ws.Range("E31").Value = "  +8.03%  "
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").Value = "  +5.75%  "
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.55"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.28"
$ws.Range("E38").Value = "  +4.70%  "
$ws.Range("E39").Value = "  +7.25%  "
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +5.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.64"
$ws.Range("E42").Value = "  +9.09%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  +4.18%  "
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.06"
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.34"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D51").Value = "0.0₆0263"
$ws.Range("E51").Value = "  +1.69%  "
# Row 49/50: swap Optimism and ARBITRUM entries with updated values
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.536"
$ws.Range("E49").Value = "  +2.84%  "

$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("E50").Value = "  +5.45%  "
